$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range keeps its existing (unstyled) text formatting
# while preventing Excel from auto-converting numeric-looking strings
# (e.g. "44.45", "0.9979") into actual numbers.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.853.02'
$ws.Range("E2").Value = '  -3.43%  '
$ws.Range("D3").Value = '1.824.09'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '276.95'
$ws.Range("E5").Value = '  -7.67%  '
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '0.5092'
$ws.Range("E7").Value = '  -4.45%  '
$ws.Range("D8").Value = '0.3475'
$ws.Range("E8").Value = '  -6.91%  '
$ws.Range("D9").Value = '44.55'
$ws.Range("E9").Value = '  -1.67%  '
$ws.Range("D10").Value = '0.06776'
$ws.Range("E10").Value = '  -5.09%  '
$ws.Range("D11").Value = '19.78'
$ws.Range("E11").Value = '  -8.11%  '
$ws.Range("D12").Value = '0.8069'
$ws.Range("E12").Value = '  -8.93%  '
$ws.Range("D13").Value = '0.07806'
$ws.Range("E13").Value = '  -3.98%  '
$ws.Range("D14").Value = '1.812.24'
$ws.Range("E14").Value = '  -4.79%  '
$ws.Range("D15").Value = '5.064'
$ws.Range("E15").Value = '  -4.20%  '
$ws.Range("D16").Value = '87.73'
$ws.Range("E16").Value = '  -4.98%  '
$ws.Range("D17").Value = '0.9984'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '14.12'
$ws.Range("E18").Value = '  -4.93%  '
$ws.Range("D19").Value = '0.000008042'
$ws.Range("E19").Value = '  -5.21%  '
$ws.Range("D20").Value = '0.9963'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '25.902.69'
$ws.Range("E21").Value = '  -3.37%  '
$ws.Range("D22").Value = '4.755'
$ws.Range("E22").Value = '  -4.28%  '
$ws.Range("D23").Value = '9.989'
$ws.Range("E23").Value = '  -6.01%  '
$ws.Range("D24").Value = '6.177'
$ws.Range("E24").Value = '  -3.17%  '
$ws.Range("D25").Value = '2.346'
$ws.Range("E25").Value = '  +2.98%  '
$ws.Range("D26").Value = '142.66'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").Value = '1.663'
$ws.Range("E27").Value = '  -4.39%  '
$ws.Range("D28").Value = '17.14'
$ws.Range("E28").Value = '  -4.57%  '
$ws.Range("D29").Value = '109.26'
$ws.Range("E29").Value = '  -3.82%  '
$ws.Range("D30").Value = '4.309'
$ws.Range("E30").Value = '  -8.17%  '
$ws.Range("D31").Value = '4.279'
$ws.Range("E31").Value = '  -7.36%  '
$ws.Range("D32").Value = '0.08750'
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").Value = '0.04851'
$ws.Range("E33").Value = '  -3.41%  '
$ws.Range("D34").Value = '1.162'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("D35").Value = '0.7249'
$ws.Range("E35").Value = '  -10.16%  '
$ws.Range("D36").Value = '2.856'
$ws.Range("E36").Value = '  -3.14%  '
$ws.Range("D37").Value = '3.171'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '0.9959'
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.388'
$ws.Range("E39").Value = '  -11.06%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01845'
$ws.Range("E40").Value = '  -4.87%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5092'
$ws.Range("E41").Value = '  -16.40%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '0.9443'
$ws.Range("E42").Value = '  -11.37%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '116.41'
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.191'
$ws.Range("E44").Value = '  -4.36%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '7.958'
$ws.Range("E45").Value = '  -9.06%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '0.9967'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").Value = '0.1361'
$ws.Range("E47").Value = '  -8.46%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4473'
$ws.Range("E48").Value = '  -15.17%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.345'
$ws.Range("E49").Value = '  -6.18%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '36.13'
$ws.Range("E50").Value = '  -3.02%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05915'
$ws.Range("E51").Value = '  -2.42%  '

# Restore default (unstyled) cell style now that values are committed as text
$dataRange.Style = "Normal"
